$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 39017.953
$ws.Range("I98").Value = 1750
$ws.Range("J98").Value = 99578.375
$ws.Range("K98").Value = 1750
$ws.Range("L98").Value = 99578.375
$ws.Range("M98").Value = -252
$ws.Range("N98").Value = -102574.375
$ws.Range("H107").Value = 4145.1035
$ws.Range("I107").Value = 4431.077
$ws.Range("J107").Value = 1666.6666
$ws.Range("K107").Value = 4431.077
$ws.Range("L107").Value = 1666.6666
$ws.Range("M107").Value = -2511.077
$ws.Range("N107").Value = -5506.6666
$ws.Range("H122").Value = 39017.953
$ws.Range("I122").Value = 1750
$ws.Range("J122").Value = 99578.375
$ws.Range("K122").Value = 5250
$ws.Range("L122").Value = 298735.125
$ws.Range("M122").Value = -2800
$ws.Range("N122").Value = -303635.125
$ws.Range("H137").Value = 6188.531
$ws.Range("I137").Value = 5375.484
$ws.Range("J137").Value = 7588.778
$ws.Range("K137").Value = 16126.452
$ws.Range("L137").Value = 22766.334
$ws.Range("M137").Value = -13576.452
$ws.Range("N137").Value = -27866.334
$ws.Range("H138").Value = 2055.041
$ws.Range("I138").Value = 1936.15
$ws.Range("J138").Value = 2099.9058
$ws.Range("K138").Value = 5808.450000000001
$ws.Range("L138").Value = 6299.7174
$ws.Range("M138").Value = -668.4500000000007
$ws.Range("N138").Value = -16579.7174

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14111.138
$ws.Range("I32").Value = 13014.667
$ws.Range("J32").Value = 22334.666
$ws.Range("K32").Value = 13014.667
$ws.Range("L32").Value = 22334.666
$ws.Range("M32").Value = -12727.667
$ws.Range("N32").Value = -22908.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2793.375
$ws.Range("I134").Value = 2387.3333
$ws.Range("J134").Value = 3636.6924
$ws.Range("K134").Value = 7161.999899999999
$ws.Range("L134").Value = 10910.0772
$ws.Range("M134").Value = -4626.999899999999
$ws.Range("N134").Value = -15980.0772

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2625.7778
$ws.Range("I22").Value = 386.4
$ws.Range("J22").Value = 5425
$ws.Range("K22").Value = 386.4
$ws.Range("L22").Value = 5425
$ws.Range("M22").Value = -36.39999999999998
$ws.Range("N22").Value = -6125
$ws.Range("H31").Value = 3971980.8
$ws.Range("I31").Value = 2241.8
$ws.Range("J31").Value = 5654073.5
$ws.Range("K31").Value = 2241.8
$ws.Range("L31").Value = 5654073.5
$ws.Range("M31").Value = -1946.8
$ws.Range("N31").Value = -5654663.5
$ws.Range("H34").Value = 3971980.8
$ws.Range("I34").Value = 2241.8
$ws.Range("J34").Value = 5654073.5
$ws.Range("K34").Value = 2241.8
$ws.Range("L34").Value = 5654073.5
$ws.Range("M34").Value = -2039.8
$ws.Range("N34").Value = -5654477.5
$ws.Range("H41").Value = 23959.6
$ws.Range("J41").Value = 28934.75
$ws.Range("L41").Value = 28934.75
$ws.Range("N41").Value = -29790.75
$ws.Range("H122").Value = 61180
$ws.Range("I122").Value = 80945.734
$ws.Range("J122").Value = 1882.8
$ws.Range("K122").Value = 242837.202
$ws.Range("L122").Value = 5648.4
$ws.Range("M122").Value = -240387.202
$ws.Range("N122").Value = -10548.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5945.423
$ws.Range("I3").Value = 2036.6666
$ws.Range("J3").Value = 8014.7646
$ws.Range("K3").Value = 6109.9998
$ws.Range("L3").Value = 24044.2938
$ws.Range("M3").Value = -5997.9998
$ws.Range("N3").Value = -24268.2938
$ws.Range("H5").Value = 4411.8076
$ws.Range("J5").Value = 661.1111
$ws.Range("L5").Value = 1983.3333
$ws.Range("N5").Value = -2207.3333
$ws.Range("H33").Value = 19067952
$ws.Range("I33").Value = 30
$ws.Range("J33").Value = 22245940
$ws.Range("K33").Value = 180
$ws.Range("L33").Value = 133475640
$ws.Range("M33").Value = 103
$ws.Range("N33").Value = -133476206
$ws.Range("H64").Value = 4750.5
$ws.Range("I64").Value = 2700
$ws.Range("J64").Value = 5889.6665
$ws.Range("K64").Value = 8100
$ws.Range("L64").Value = 17668.9995
$ws.Range("M64").Value = -7830
$ws.Range("N64").Value = -18208.9995
$ws.Range("H67").Value = 4750.5
$ws.Range("I67").Value = 2700
$ws.Range("J67").Value = 5889.6665
$ws.Range("K67").Value = 8100
$ws.Range("L67").Value = 17668.9995
$ws.Range("M67").Value = -7164
$ws.Range("N67").Value = -19540.9995
$ws.Range("H81").Value = 2153.3333
$ws.Range("I81").Value = 566.6667
$ws.Range("J81").Value = 2946.6667
$ws.Range("K81").Value = 1700.0001
$ws.Range("L81").Value = 8840.000100000001
$ws.Range("M81").Value = -577.0001
$ws.Range("N81").Value = -11086.0001
$ws.Range("H84").Value = 2153.3333
$ws.Range("I84").Value = 566.6667
$ws.Range("J84").Value = 2946.6667
$ws.Range("K84").Value = 5100.0003
$ws.Range("L84").Value = 26520.0003
$ws.Range("M84").Value = 515.9997000000003
$ws.Range("N84").Value = -37752.0003
$ws.Range("H109").Value = 2578.0588
$ws.Range("I109").Value = 1503.3846
$ws.Range("J109").Value = 3243.3333
$ws.Range("K109").Value = 4510.1538
$ws.Range("L109").Value = 9729.999899999999
$ws.Range("M109").Value = -3470.1538
$ws.Range("N109").Value = -11809.9999
$ws.Range("H113").Value = 7161.875
$ws.Range("I113").Value = 13097.375
$ws.Range("J113").Value = 1226.375
$ws.Range("K113").Value = 39292.125
$ws.Range("L113").Value = 3679.125
$ws.Range("M113").Value = -37122.125
$ws.Range("N113").Value = -8019.125
$ws.Range("H114").Value = 1603.7037
$ws.Range("I114").Value = 1205.5834
$ws.Range("J114").Value = 1922.2
$ws.Range("K114").Value = 3616.7502
$ws.Range("L114").Value = 5766.6
$ws.Range("M114").Value = -362.7501999999999
$ws.Range("N114").Value = -12274.6
$ws.Range("H117").Value = 2441.125
$ws.Range("J117").Value = 2535.2666
$ws.Range("L117").Value = 7605.7998
$ws.Range("N117").Value = -14489.7998
$ws.Range("H122").Value = 3022.9556
$ws.Range("I122").Value = 643.30304
$ws.Range("J122").Value = 9567
$ws.Range("K122").Value = 5789.72736
$ws.Range("L122").Value = 86103
$ws.Range("M122").Value = -3339.72736
$ws.Range("N122").Value = -91003
$ws.Range("H132").Value = 2516
$ws.Range("I132").Value = 1314.2858
$ws.Range("J132").Value = 5320
$ws.Range("K132").Value = 11828.5722
$ws.Range("L132").Value = 47880
$ws.Range("M132").Value = -9298.5722
$ws.Range("N132").Value = -52940
$ws.Range("H133").Value = 5218.3335
$ws.Range("I133").Value = 5162
$ws.Range("J133").Value = 5500
$ws.Range("K133").Value = 15486
$ws.Range("L133").Value = 16500
$ws.Range("M133").Value = -10426
$ws.Range("N133").Value = -26620
$ws.Range("H135").Value = 4411.8076
$ws.Range("J135").Value = 661.1111
$ws.Range("L135").Value = 5949.9999
$ws.Range("N135").Value = -11019.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5447.8887
$ws.Range("I80").Value = 4071.5
$ws.Range("K80").Value = 4071.5
$ws.Range("M80").Value = -3073.5
$ws.Range("H83").Value = 5447.8887
$ws.Range("I83").Value = 4071.5
$ws.Range("K83").Value = 20357.5
$ws.Range("M83").Value = -15365.5
$ws.Range("H126").Value = 16024.75
$ws.Range("I126").Value = 37341.332
$ws.Range("K126").Value = 112023.996
$ws.Range("M126").Value = -109553.996
